$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two claim numbers so FT validation accepts both "daños" and "robo" claims
$ws.Range("B5").Value = "'0420172008479   "
$ws.Range("B6").Value = "'1120170200933"

# Move active selection to B8 (was C8)
$ws.Range("B8").Select()
